# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextCell "D2" '42.453.33'
Set-TextCell "E2" '  -1.78%  '
Set-TextCell "D3" '2.185.06'
Set-TextCell "E3" '  -2.46%  '
Set-TextCell "E4" '  -0.07%  '
Set-TextCell "D5" '251.84'
Set-TextCell "E5" '  +2.20%  '
Set-TextCell "D6" '0.606'
Set-TextCell "E6" '  -2.34%  '
Set-TextCell "D7" '75.09'
Set-TextCell "E7" '  -1.13%  '
Set-TextCell "E8" '  -0.01%  '
Set-TextCell "E9" '  -5.82%  '
Set-TextCell "D10" '40.23'
Set-TextCell "E10" '  -3.29%  '
Set-TextCell "D11" '0.0910'
Set-TextCell "E11" '  -2.98%  '
Set-TextCell "E12" '  -0.47%  '
Set-TextCell "D13" '6.76'
Set-TextCell "E13" '  -3.40%  '
Set-TextCell "D14" '2.513.11'
Set-TextCell "E14" '  -2.41%  '
Set-TextCell "D15" '14.17'
Set-TextCell "E15" '  -4.11%  '
Set-TextCell "D16" '2.182.24'
Set-TextCell "E16" '  -2.33%  '
Set-TextCell "D17" '0.769'
Set-TextCell "E17" '  -5.65%  '
Set-TextCell "D18" '42.377.65'
Set-TextCell "E18" '  -1.69%  '
Set-TextCell "E19" '  -3.63%  '
Set-TextCell "D20" '70.82'
Set-TextCell "E20" '  -0.47%  '
Set-TextCell "D21" '5.85'
Set-TextCell "E21" '  -2.77%  '
Set-TextCell "D22" '226.61'
Set-TextCell "E22" '  -1.93%  '
Set-TextCell "D23" '9.37'
Set-TextCell "E23" '  -11.57%  '
Set-TextCell "D24" '2.11'
Set-TextCell "E24" '  -2.78%  '
Set-TextCell "E25" '  -0.01%  '
Set-TextCell "D26" '10.44'
Set-TextCell "E26" '  -4.98%  '
Set-TextCell "D27" '3.42'
Set-TextCell "E27" '  +2.09%  '
Set-TextCell "E28" '  -4.21%  '
Set-TextCell "B29" 'Monero'
Set-TextCell "C29" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell "D29" '171.93'
Set-TextCell "E29" '  -1.46%  '
Set-TextCell "B30" 'Toncoin'
Set-TextCell "C30" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell "D30" '2.14'
Set-TextCell "E30" '  -2.53%  '
Set-TextCell "B31" 'InjectiveProtocol'
Set-TextCell "C31" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell "D31" '36.84'
Set-TextCell "E31" '  -1.62%  '
Set-TextCell "D32" '20.02'
Set-TextCell "E32" '  -1.82%  '
Set-TextCell "D33" '0.0818'
Set-TextCell "E33" '  +2.62%  '
Set-TextCell "D34" '5.13'
Set-TextCell "E34" '  -5.00%  '
Set-TextCell "E35" '  -2.13%  '
Set-TextCell "E36" '  -4.42%  '
Set-TextCell "D37" '4.20'
Set-TextCell "E37" '  -3.25%  '
Set-TextCell "E38" '  +0.21%  '
Set-TextCell "D39" '12.00'
Set-TextCell "E39" '  -9.38%  '
Set-TextCell "D40" '2.06'
Set-TextCell "E40" '  -3.87%  '
Set-TextCell "E41" '  +10.30%  '
Set-TextCell "E42" '  -8.03%  '
Set-TextCell "E43" '  -3.26%  '
Set-TextCell "E44" '  -3.17%  '
Set-TextCell "D45" '101.52'
Set-TextCell "E45" '  -4.01%  '
Set-TextCell "E46" '  -2.67%  '
Set-TextCell "E47" '  -4.37%  '
Set-TextCell "D48" '0.455'
Set-TextCell "E48" '  -0.30%  '
Set-TextCell "D49" '1.09'
Set-TextCell "E49" '  -2.28%  '
Set-TextCell "D50" '1.12'
Set-TextCell "E50" '  -2.46%  '
Set-TextCell "E51" '  -0.77%  '

Write-Host "Applied cryptos update"
